$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing weekly figures (columns G = Waargenomen, H = Verwacht) ---
# The I column holds shared formulas (=G-H) that recalculate automatically.
$ws.Range("G3").Value = 3217
$ws.Range("G4").Value = 3615
$ws.Range("G5").Value = 4459
$ws.Range("H7").Value = 2909
$ws.Range("H8").Value = 3010
$ws.Range("G9").Value = 3906
$ws.Range("G11").Value = 2984
$ws.Range("H11").Value = 2933
$ws.Range("H12").Value = 3050
$ws.Range("G13").Value = 2770
$ws.Range("G14").Value = 2725
$ws.Range("G17").Value = 2692
$ws.Range("G19").Value = 2636
$ws.Range("G20").Value = 2614
$ws.Range("H20").Value = 2856
$ws.Range("G21").Value = 2526
$ws.Range("G22").Value = 2670
$ws.Range("G23").Value = 2657
$ws.Range("G24").Value = 2634
$ws.Range("G25").Value = 3202
$ws.Range("G26").Value = 2836
$ws.Range("G27").Value = 2715
$ws.Range("G28").Value = 2661

# --- Insert a new row for week 37 right after the current last data row (28) ---
# This pushes the totals row (old row 30, with a blank row 29 gap) down to row 31.
$ws.Rows("29").Insert()

$ws.Range("F29").Value = 37
$ws.Range("G29").Value = 2704
$ws.Range("H29").Value = 2844
$ws.Range("I29").Formula = "=G29-H29"

# --- Move the active selection to I13, matching the author's final cursor position ---
$ws.Range("I13").Select()
